$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1 (API-Testing): add a new "PetGetTest" / "Testing Empty" test row
# ---------------------------------------------------------------------------
$ws1.Range("A5").Value = "PetGetTest"
$ws1.Range("B5").Value = "Testing Empty"
$ws1.Range("C5").Value = "https://live.virtualandemo.com/api/pets/findByTags?tags=light-grey"
$ws1.Range("D5").Value = "application/json"
$ws1.Range("F5").Value = "get_response_empty.json"
$ws1.Range("H5").Value = "GET"
$ws1.Range("J5").Value = 200

# Hyperlink the URL cell, re-using the same target address as the other rows
$ws1.Hyperlinks.Add($ws1.Range("C5"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey") | Out-Null

# Match the formatting used by the other rows in these columns (font, fill...)
$ws1.Range("C3").Copy() | Out-Null
$ws1.Range("C5").PasteSpecial(-4122) | Out-Null
$ws1.Range("D3").Copy() | Out-Null
$ws1.Range("D5").PasteSpecial(-4122) | Out-Null
$ws1.Range("F3").Copy() | Out-Null
$ws1.Range("F5").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2 (API-Testing-Sheet2-Duplicate): add the matching "PetEmptyTest" row
# ---------------------------------------------------------------------------
$ws2.Range("A5").Value = "PetEmptyTest"
$ws2.Range("B5").Value = "API get testing"
$ws2.Range("C5").Value = "https://live.virtualandemo.com/api/pets/findByTags?tags=light-grey"
$ws2.Range("D5").Value = "application/json"
$ws2.Range("F5").Value = "get_response_empty.json"
$ws2.Range("H5").Value = "GET"
$ws2.Range("J5").Value = 200

$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey") | Out-Null

$ws2.Range("C3").Copy() | Out-Null
$ws2.Range("C5").PasteSpecial(-4122) | Out-Null
$ws2.Range("D3").Copy() | Out-Null
$ws2.Range("D5").PasteSpecial(-4122) | Out-Null
$ws2.Range("F3").Copy() | Out-Null
$ws2.Range("F5").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("M5").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("C13").Select() | Out-Null

$ws1.Activate() | Out-Null
